$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Remove the obsolete "Flexible" label from the "Tipo_Horario" column (N) for
# rows 62-69; replace with "Frecuencia", move the existing duration value
# from the "Duracion_Trayecto_Min" column (O) into the new
# "Frecuencia_Min" column (P), and reset O to 0 -- matching the pattern
# already used by rows 60-61.
$frecuenciaMin = @{
    62 = 0.010416666666666666
    63 = 0.010416666666666666
    64 = 0.008333333333333333
    65 = 0.008333333333333333
    66 = 0.017361111111111112
    67 = 0.017361111111111112
    68 = 0.024305555555555556
    69 = 0.024305555555555556
}

for ($r = 62; $r -le 69; $r++) {
    $ws.Cells.Item($r, 14).Value = "Frecuencia"
    $ws.Cells.Item($r, 15).Value = 0
    $ws.Cells.Item($r, 16).Value2 = $frecuenciaMin[$r]
    # The new Frecuencia_Min cell should carry the same time format as the
    # Duracion_Trayecto_Min cell next to it.
    $ws.Cells.Item($r, 15).Copy()
    $ws.Cells.Item($r, 16).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# Widen column P (Frecuencia_Min) now that it holds real data for more rows.
$ws.Columns.Item(16).ColumnWidth = 24.5703125

# Update the sheet view/selection to reflect where the user was working.
$ws.Range("O63").Select()
$excel.ActiveWindow.ScrollRow = 43
